$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet: rows 2-44 (row 45 & 46 removed) ---
$recoData = @(
    @("BRVM - SERVICES PUBLICS", 0, 6, 2424.84, 99.66, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 3, 1915.69, 637.7, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 3, 1790, 595, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI", 0, 3, 1735, 580, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI", 0, 3, 1655, 580, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 3, 1645, 545, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 3, 1104.05, 371.75, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 3, 1046.41, 353.7, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 3, 972.82, 326.77, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE", 0, 3, 789.4, 262.08, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE", 0, 3, 652.31, 217.18, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI", 0, 1, 635, 635, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL", 0, 3, 569.28, 189.81, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 3, 410.05, 137.18, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 3, 387.59, 129.75, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 3, 366.14, 122.1, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 3, 359.83, 119.99, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 3, 329.67, 110.63, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 3, 317, 106.5, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 3, 277.39, 93.32, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 2, 1, 7.04, -5.88, "🟡 Observer", "👀 À surveiller"),
    @("VIVO ENERGY CI (SHEC)", 1, 0, 5.73, 5.73, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 1, 0, 5.45, 5.45, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 1, 1, 3.96, -3.45, "🟡 Observer", "👀 À surveiller"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 1, 0, 3.54, 3.54, "🟡 Observer", "➖ Neutre"),
    @("PALM CI (PALC)", 1, 0, 3.41, 3.41, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 1, 0, 3.33, 3.33, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 3.09, 3.09, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 1, 0, 3.05, 3.05, "🟡 Observer", "➖ Neutre"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 1, 0, 2.93, 2.93, "🟡 Observer", "➖ Neutre"),
    @("SICABLE CI (CABC)", 1, 0, 2.69, 2.69, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 1, 0, 1.96, 1.96, "🟡 Observer", "➖ Neutre"),
    @("SAPH CI (SPHC)", 1, 1, 1.34, 4.03, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 3, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("BICI CI (BICC)", 0, 1, -0.6, -0.6, "🟡 Observer", "➖ Neutre"),
    @("ORAGROUP TOGO (ORGT)", 0, 1, -0.6, -0.6, "🟡 Observer", "➖ Neutre"),
    @("SERVAIR ABIDJAN CI (ABJC)", 0, 1, -1.42, -1.42, "🟡 Observer", "➖ Neutre"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -1.75, -1.75, "🟡 Observer", "➖ Neutre"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 0, 1, -1.87, -1.87, "🟡 Observer", "➖ Neutre"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 0, 1, -2, -2, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 0, 1, -2.02, -2.02, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 0, 2, -3.94, -1.97, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 1, 2, -11.22, -7.33, "🟡 Observer", "👀 À surveiller")
)

for ($i = 0; $i -lt $recoData.Count; $i++) {
    $r = $i + 2
    $row = $recoData[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsReco.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Remove the now-obsolete trailing rows (old dimension was A1:G46)
$wsReco.Rows("45:46").Delete()

# --- Top_YTD sheet: rows 2-11, column A (Titre) + B (Progression YTD %) ---
$ytdData = @(
    @("BRVM - SERVICES PUBLICS", 422362.11),
    @("BRVM - AUTRES SECTEURS", 40186.58),
    @("NEI-CEDA CI", 33711.75),
    @("UNIWAX CI", 31112),
    @("SETAO CI", 27535.2),
    @("AIR LIQUIDE CI", 27151.25),
    @("BRVM - DISTRIBUTION", 10150.6),
    @("BRVM - TRANSPORT", 8938.92),
    @("BRVM - AGRICULTURE", 7536.75),
    @("BRVM - INDUSTRIE", 4688.46)
)

for ($i = 0; $i -lt $ytdData.Count; $i++) {
    $r = $i + 2
    $row = $ytdData[$i]
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
}
